$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 9618215
$ws.Range("I41").Value = 13892522
$ws.Range("J41").Value = 1025
$ws.Range("K41").Value = 13892522
$ws.Range("L41").Value = 1025
$ws.Range("M41").Value = -13892082
$ws.Range("N41").Value = -1905
$ws.Range("H98").Value = 62503130
$ws.Range("I98").Value = 76926696
$ws.Range("K98").Value = 76926696
$ws.Range("M98").Value = -76925198
$ws.Range("H113").Value = 150007740
$ws.Range("I113").Value = 2805
$ws.Range("J113").Value = 214295570
$ws.Range("K113").Value = 2805
$ws.Range("L113").Value = 214295570
$ws.Range("M113").Value = 449
$ws.Range("N113").Value = -214302078
$ws.Range("H122").Value = 62503130
$ws.Range("I122").Value = 76926696
$ws.Range("K122").Value = 230780088
$ws.Range("M122").Value = -230777638
$ws.Range("H132").Value = 1283.3541
$ws.Range("I132").Value = 1233.4524
$ws.Range("K132").Value = 3700.357199999999
$ws.Range("M132").Value = -1170.357199999999
$ws.Range("H135").Value = 303629.34
$ws.Range("I135").Value = 323192.53
$ws.Range("K135").Value = 2908732.77
$ws.Range("M135").Value = -2906197.77
$ws.Range("H137").Value = 5083.5
$ws.Range("I137").Value = 13001
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 39003
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -36453
$ws.Range("N137").Value = -15600
$ws.Range("H138").Value = 4172.323
$ws.Range("J138").Value = 7617.7417
$ws.Range("L138").Value = 22853.2251
$ws.Range("N138").Value = -33133.2251

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1589.79
$ws.Range("I32").Value = 1573.9791
$ws.Range("K32").Value = 1573.9791
$ws.Range("M32").Value = -1286.9791
$ws.Range("H45").Value = 7545.273
$ws.Range("I45").Value = 6333.1113
$ws.Range("J45").Value = 13000
$ws.Range("K45").Value = 6333.1113
$ws.Range("L45").Value = 13000
$ws.Range("M45").Value = -5956.1113
$ws.Range("N45").Value = -13754
$ws.Range("H61").Value = 6686.4653
$ws.Range("I61").Value = 3754.1035
$ws.Range("K61").Value = 3754.1035
$ws.Range("M61").Value = -3542.1035
$ws.Range("H63").Value = 1871.6666
$ws.Range("I63").Value = 1246
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 1246
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -560
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 1871.6666
$ws.Range("I66").Value = 1246
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 6230
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -2798
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 25941.707
$ws.Range("I74").Value = 32109.25
$ws.Range("K74").Value = 32109.25
$ws.Range("M74").Value = -31235.25
$ws.Range("H77").Value = 25941.707
$ws.Range("I77").Value = 32109.25
$ws.Range("K77").Value = 160546.25
$ws.Range("M77").Value = -156178.25
$ws.Range("H97").Value = 16667137
$ws.Range("I97").Value = 575
$ws.Range("J97").Value = 27778178
$ws.Range("K97").Value = 575
$ws.Range("L97").Value = 27778178
$ws.Range("M97").Value = -79
$ws.Range("N97").Value = -27779170
$ws.Range("H122").Value = 4174.8
$ws.Range("I122").Value = 2438.5454
$ws.Range("K122").Value = 7315.6362
$ws.Range("M122").Value = -4865.6362
$ws.Range("H132").Value = 10246.154
$ws.Range("I132").Value = 11835.777
$ws.Range("K132").Value = 35507.331
$ws.Range("M132").Value = -32977.331
$ws.Range("H136").Value = 6686.4653
$ws.Range("I136").Value = 3754.1035
$ws.Range("K136").Value = 11262.3105
$ws.Range("M136").Value = -8712.3105

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 297.5
$ws.Range("I22").Value = 297.5
$ws.Range("K22").Value = 297.5
$ws.Range("M22").Value = -124.5
$ws.Range("H134").Value = 6518.2617
$ws.Range("I134").Value = 2415.6428
$ws.Range("K134").Value = 7246.928400000001
$ws.Range("M134").Value = -4711.928400000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5673.4307
$ws.Range("I31").Value = 2403.1086
$ws.Range("K31").Value = 2403.1086
$ws.Range("M31").Value = -2108.1086
$ws.Range("H34").Value = 5673.4307
$ws.Range("I34").Value = 2403.1086
$ws.Range("K34").Value = 2403.1086
$ws.Range("M34").Value = -2201.1086
$ws.Range("H58").Value = 8201372
$ws.Range("I58").Value = 12501705
$ws.Range("K58").Value = 12501705
$ws.Range("M58").Value = -12501502
$ws.Range("H99").Value = 7311.2
$ws.Range("I99").Value = 3699.3333
$ws.Range("K99").Value = 3699.3333
$ws.Range("M99").Value = -2201.3333
$ws.Range("H107").Value = 2253.4443
$ws.Range("I107").Value = 561.8570999999999
$ws.Range("K107").Value = 561.8570999999999
$ws.Range("M107").Value = 1358.1429
$ws.Range("H126").Value = 7311.2
$ws.Range("I126").Value = 3699.3333
$ws.Range("K126").Value = 11097.9999
$ws.Range("M126").Value = -8627.999899999999
$ws.Range("H134").Value = 4733.1694
$ws.Range("I134").Value = 1941.8108
$ws.Range("K134").Value = 5825.4324
$ws.Range("M134").Value = -3290.4324
$ws.Range("H136").Value = 8201372
$ws.Range("I136").Value = 12501705
$ws.Range("K136").Value = 37505115
$ws.Range("M136").Value = -37502565

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 12821380
$ws.Range("I14").Value = 12821380
$ws.Range("K14").Value = 38464140
$ws.Range("M14").Value = -38463967
$ws.Range("H54").Value = 1406
$ws.Range("H140").Value = 79620.38
$ws.Range("I140").Value = 112341
$ws.Range("K140").Value = 337023
$ws.Range("M140").Value = -331843

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 45908.91
$ws.Range("I52").Value = 32142.715
$ws.Range("K52").Value = 32142.715
$ws.Range("M52").Value = -31883.715
$ws.Range("H97").Value = 1589.3182
$ws.Range("I97").Value = 1577.1578
$ws.Range("K97").Value = 1577.1578
$ws.Range("M97").Value = -1081.1578
$ws.Range("H132").Value = 9097.857
$ws.Range("I132").Value = 3428.889
$ws.Range("K132").Value = 10286.667
$ws.Range("M132").Value = -7756.667000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7858.5454
$ws.Range("I7").Value = 5611
$ws.Range("K7").Value = 5611
$ws.Range("M7").Value = -5499
$ws.Range("H40").Value = 6007.316
$ws.Range("I40").Value = 3340.25
$ws.Range("K40").Value = 3340.25
$ws.Range("M40").Value = -3204.25
$ws.Range("H100").Value = 4166.625
$ws.Range("I100").Value = 3409.4
$ws.Range("K100").Value = 3409.4
$ws.Range("M100").Value = -2868.4
$ws.Range("H126").Value = 7858.5454
$ws.Range("I126").Value = 5611
$ws.Range("K126").Value = 16833
$ws.Range("M126").Value = -14363

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 88090.664
$ws.Range("I62").Value = 129636
$ws.Range("K62").Value = 129636
$ws.Range("M62").Value = -129012
$ws.Range("H65").Value = 88090.664
$ws.Range("I65").Value = 129636
$ws.Range("K65").Value = 648180
$ws.Range("M65").Value = -645060
$ws.Range("H107").Value = 16667649
$ws.Range("I107").Value = 976.0909
$ws.Range("K107").Value = 2928.2727
$ws.Range("M107").Value = -1008.2727

Write-Host "Applied all market-price updates."